$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells keep their original text formatting
# (values like "116.95" would otherwise be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '52.051.73'
$ws.Range("E2").Value = '  +5.23%  '
$ws.Range("D3").Value = '2.787.97'
$ws.Range("E3").Value = '  +6.25%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '116.95'
$ws.Range("E5").Value = '  +4.97%  '
$ws.Range("D6").Value = '339.77'
$ws.Range("E6").Value = '  +4.64%  '
$ws.Range("D7").Value = '0.542'
$ws.Range("E7").Value = '  +3.76%  '
$ws.Range("D9").Value = '0.578'
$ws.Range("E9").Value = '  +6.65%  '
$ws.Range("D10").Value = '42.49'
$ws.Range("E10").Value = '  +8.00%  '
$ws.Range("D11").Value = '0.0872'
$ws.Range("E11").Value = '  +8.16%  '
$ws.Range("D12").Value = '20.06'
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("E13").Value = '  +2.44%  '
$ws.Range("D14").Value = '7.65'
$ws.Range("E14").Value = '  +4.48%  '
$ws.Range("D15").Value = '3.218.36'
$ws.Range("E15").Value = '  +6.07%  '
$ws.Range("D16").Value = '2.815.27'
$ws.Range("E16").Value = '  +7.07%  '
$ws.Range("D17").Value = '0.888'
$ws.Range("E17").Value = '  +4.69%  '
$ws.Range("D18").Value = '51.853.28'
$ws.Range("E18").Value = '  +5.04%  '
$ws.Range("E19").Value = '  +12.30%  '
$ws.Range("D20").Value = '13.45'
$ws.Range("E20").Value = '  +3.96%  '
$ws.Range("D21").Value = '6.96'
$ws.Range("E21").Value = '  +4.55%  '
$ws.Range("D22").Value = '0.0₃0981'
$ws.Range("E22").Value = '  +3.92%  '
$ws.Range("D23").Value = '277.96'
$ws.Range("E23").Value = '  +4.37%  '
$ws.Range("D24").Value = '70.28'
$ws.Range("E24").Value = '  +2.42%  '
$ws.Range("D25").Value = '2.77'
$ws.Range("E25").Value = '  +9.57%  '
$ws.Range("D26").Value = '26.90'
$ws.Range("E26").Value = '  +3.65%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '10.22'
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("E29").Value = '  +1.16%  '
$ws.Range("E30").Value = '  +3.53%  '
$ws.Range("D31").Value = '35.19'
$ws.Range("E31").Value = '  +2.24%  '
$ws.Range("D32").Value = '50.23'
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").Value = '5.66'
$ws.Range("E33").Value = '  +3.53%  '
$ws.Range("D34").Value = '0.0822'
$ws.Range("E34").Value = '  +2.17%  '
$ws.Range("D35").Value = '2.14'
$ws.Range("E35").Value = '  +5.64%  '
$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").Value = '19.09'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("D38").Value = '3.30'
$ws.Range("E38").Value = '  +7.39%  '
$ws.Range("D39").Value = '4.99'
$ws.Range("E39").Value = '  +1.50%  '
$ws.Range("D40").Value = '2.72'
$ws.Range("E40").Value = '  +26.18%  '
$ws.Range("D41").Value = '0.0368'
$ws.Range("E41").Value = '  +14.01%  '
$ws.Range("D42").Value = '23.68'
$ws.Range("E42").Value = '  +5.18%  '
$ws.Range("E43").Value = '  +6.20%  '
$ws.Range("D44").Value = '127.02'
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("E45").Value = '  +3.47%  '
$ws.Range("D46").Value = '2.106.13'
$ws.Range("E46").Value = '  +3.42%  '
$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  +4.46%  '
$ws.Range("E48").Value = '  +3.43%  '
$ws.Range("D49").Value = '5.57'
$ws.Range("E49").Value = '  +7.46%  '
$ws.Range("D50").Value = '0.912'
$ws.Range("E50").Value = '  +21.73%  '
$ws.Range("D51").Value = '8.91'
$ws.Range("E51").Value = '  +0.96%  '
